# "First cut with expander code" - add a new "With Expander" column (I) to
# the Wemos pin-mapping sheet, documenting which board pin should be used
# for each expander signal.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell, bold + underlined like the other section headers.
$ws.Range("I1").Value = "With Expander"
$ws.Range("I1").Font.Bold = $true
$ws.Range("I1").Font.Underline = $true

# New-string cells first (so the shared-string table grows in the same
# order the original author's Excel session produced them), then the
# cell that reuses the pre-existing "Sound-Rx" string.
$ws.Range("I7").Value = "SDA"
$ws.Range("I8").Value = "SCL"
$ws.Range("I5").Value = "Sound-Tx (not used)"
$ws.Range("I9").Value = "INT"
$ws.Range("I6").Value = "Sound-Rx"

# Leave the selection where the author's session left it.
$ws.Range("F9").Select() | Out-Null
